$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.014.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.21%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.918.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.83%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.34%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'588.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.51%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'146.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.39%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.11%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +1.00%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.917.21"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.70%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'7.04"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.96%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.153"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +7.16%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.437"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.62%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.0000239"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +6.58%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'32.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.87%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -1.44%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.400.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.33%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'61.979.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.53%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  -1.02%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'2.915.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.89%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'436.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.61%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'13.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.62%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.660"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.84%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'6.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.65%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'80.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'10.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -3.31%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'11.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.01%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -2.37%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.04%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.0000108"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +23.98%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +6.11%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'2.58"
$ws.Range("D31").Style = "Normal"
$ws.Range("E32").Value = "'  -0.09%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.110"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +3.45%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'26.05"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.34%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.16%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.977"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.74%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +8.61%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'5.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.91%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'49.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.37%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +0.68%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'8.38"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.65%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -2.54%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -0.85%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'39.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.85%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'2.701.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.46%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'134.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.14%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +0.88%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'347.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -3.70%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +0.05%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +0.36%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'22.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.68%  "
$ws.Range("E51").Style = "Normal"
